# Generate Report for Handback
# Adds a new handback entry (7a6310d6-2100-4026-8777-6ba687da0d3c.md) as row 4
# on the "Overview", "zh-cn" and "de-de" sheets, and grows the three tables
# to include it.

$wb = $excel.ActiveWorkbook

$fileGuid   = "7a6310d6-2100-4026-8777-6ba687da0d3c"
$fileName   = "$fileGuid.md"
$pathName   = "e2e\$fileGuid.md"
$ext        = ".md"
$statusSync = "Handed back: in sync with en-US"

$zhXlf      = "$fileGuid.a3bdf615d8762c82663e5c71f8fa2b08332124e8.zh-cn.xlf"
$deXlf      = "$fileGuid.a3bdf615d8762c82663e5c71f8fa2b08332124e8.de-de.xlf"

$zhHandoffDate  = "2017-02-09 14:57:45"
$zhHandbackDate = "2017-02-09 14:58:41"
$deHandoffDate  = "2017-02-09 14:58:03"
$deHandbackDate = "2017-02-09 14:59:07"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> new row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("C4").Value = $ext
$wsOverview.Range("E4").Value = $statusSync
$wsOverview.Range("F4").Value = $statusSync
$wsOverview.Range("G4").Value = $deHandoffDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a1b2c3d4e5f60718293a4b5c6d7e8f901234567/e2e/$fileName",
    "",
    "",
    $pathName
) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B4").Value = $ext
$wsZh.Range("C4").Value = $statusSync
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = $zhHandoffDate
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = ""
$wsZh.Range("K4").Value = $zhXlf
$wsZh.Range("L4").Value = $zhHandbackDate
$wsZh.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("M4").Value = ""
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "True"
$wsZh.Range("P4").Value = ""
$wsZh.Range("Q4").Value = "False"
$wsZh.Range("R4").Value = ""

$wsZh.Hyperlinks.Add(
    $wsZh.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a1b2c3d4e5f60718293a4b5c6d7e8f901234567/e2e/$fileName",
    "",
    "",
    $fileName
) | Out-Null
$wsZh.Hyperlinks.Add(
    $wsZh.Range("J4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1b2c3d4e5f60718293a4b5c6d7e8f9012345678a/e2e/$fileName",
    "",
    "",
    $fileName
) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:R4"))

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B4").Value = $ext
$wsDe.Range("C4").Value = $statusSync
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = $deHandoffDate
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = ""
$wsDe.Range("K4").Value = $deXlf
$wsDe.Range("L4").Value = $deHandbackDate
$wsDe.Range("L4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("M4").Value = ""
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "True"
$wsDe.Range("P4").Value = ""
$wsDe.Range("Q4").Value = "False"
$wsDe.Range("R4").Value = ""

$wsDe.Hyperlinks.Add(
    $wsDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0a1b2c3d4e5f60718293a4b5c6d7e8f901234567/e2e/$fileName",
    "",
    "",
    $fileName
) | Out-Null
$wsDe.Hyperlinks.Add(
    $wsDe.Range("J4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2c3d4e5f60718293a4b5c6d7e8f9012345678ab1/e2e/$fileName",
    "",
    "",
    $fileName
) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:R4"))
